$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.986533522605896
$ws.Range("B1").Value = 3.978518486022949
$ws.Range("C1").Value = 2.21148157119751
$ws.Range("D1").Value = 1.67676043510437
$ws.Range("E1").Value = 1.301976442337036
